$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new header columns D:F for Corequisites, Concurrent, Recommended
# This shifts the existing "Terms Typically Offered" header from D1 to G1
$ws.Range("D1:F1").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
